## "updated activity till excel form"
## Kings XI Punjab / Mohammed Shami sheet: update the per-innings runs/balls
## figures in rows 2-4 (columns C=runs, D=balls). The source cells store
## these numeric-looking values as TEXT, so a leading apostrophe is used to
## force text entry instead of letting Excel auto-convert them to numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: runs 0 -> 2, balls 1 -> 2
$ws.Range("C2").Value = "'2"
$ws.Range("D2").Value = "'2"

# Row 3: runs 2 -> 0, balls 2 -> 0
$ws.Range("C3").Value = "'0"
$ws.Range("D3").Value = "'0"

# Row 4: balls 0 -> 1
$ws.Range("D4").Value = "'1"
